$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_8a_Links")

# Insert a new row before row 11 (shifts existing rows 11.. down by one)
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new link entry
$ws.Cells.Item(11, 1).Value = "L_BMFSFJ_1"
$ws.Cells.Item(11, 2).Value = "Q_BMFSFJ"
$ws.Cells.Item(11, 3).Value = "Gesetze für mehr Frauen in Führungspositionen"
$ws.Cells.Item(11, 4).Value = "Laws for more women in management positions (only available in German)"
$ws.Cells.Item(11, 5).Value = "http://www.bmfsfj.de/frauen-in-fuehrungspositionen"
$ws.Cells.Item(11, 6).Value = ""
